$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Header block: replace the first paragraph ("Gustavo Camilo Ferreira")
#    with 5 new title/author paragraphs followed by the (now bold) original
#    "Aluno: Gustavo Camilo Ferreira" paragraph.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)

$headerXml = @"
<w:p $wns>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>TÓPICOS EM APRENDIZAGEM DE MÁQUINA</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">AVALIAÇÃO 2 </w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Prof. Alexandre Szabo</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">Aluno: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Gustavo Camilo Ferreira</w:t>
  </w:r>
</w:p>
"@

$p1.Range.InsertXML($headerXml)

# ---------------------------------------------------------------------------
# 2) Find the "G)" list paragraph (currently an empty list item carrying the
#    numPr for item G) and give it its answer text, then append 3 blank
#    "PargrafodaLista" paragraphs after it.
# ---------------------------------------------------------------------------
$gPara = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($p.Range.ListFormat.ListString -eq "2)" -and $txt.Trim().Length -eq 0) {
        $gPara = $p
        break
    }
}

$gXml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>G) Percebe-se pelas execuções da função "</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>main</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>" que a execução funciona com 100% de eficiência quando o k é igual a 1, e conforme o número k cresce, a eficiência do algoritmo decai. Portanto, k = 1 é a melhor valor para este determinado problema.</w:t>
  </w:r>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
  </w:pPr>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
  </w:pPr>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
  </w:pPr>
</w:p>
"@

$gPara.Range.InsertXML($gXml)
